$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Row 6: Task "Jogar o jogo para conhecer melhor o projeto" - update estimate and day6 effort
$ws.Range("D6").Value = 7
$ws.Range("J6").Value = 1

# Row 14: New task "Pesquisar como se joga o jogo"
$ws.Range("C14").Value = "Pesquisar como se joga o jogo"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 1
$ws.Range("J14").Value = 1

# Update the selected cell to N12
$ws.Range("N12").Select()
